$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$f1 = $sec.Footers.Item(1)
$shp = $f1.Range.InlineShapes.Item(1)
Write-Host "Type before:" $shp.Type
$floatShape = $shp.ConvertToShape()
Write-Host "Converted shape name:" $floatShape.Name
$floatShape.Name = "image2.png"
Write-Host "After set name:" $floatShape.Name
$newInline = $floatShape.ConvertToInlineShape()
Write-Host "back to inline, count:" $f1.Range.InlineShapes.Count
